# DDAf_2022_Tableau_annexe_Tab19.xlsx - "Add files via upload" commit
#
# Changes applied (per canonical OOXML diff):
#  1. Fix "Etats" -> "États" typos in several region-group labels
#     (column B, rows 93, 94, 97, 98).
#  2. Row 96 previously (incorrectly) repeated the row-95 label
#     ("Afrique, pays en développement sans littoral"); it is corrected
#     to the missing "RDM, pays en développement sans littoral" label.
#  3. The footnote text in A104/A105 is swapped: the "Union douanière"
#     footnote now appears before the "Responsabilité" footnote.
#  4. Row 92 (RDM, pays les moins avancés) data values C:J are updated
#     to refreshed figures.
#  5. Best-effort attempt at the saved window size in bookViews
#     (cosmetic last-view-state; engine may not persist it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & 2: region-group labels (column B) --------------------------------
$ws.Range("B93").Value = "Afrique, petits États insulaires en développement"
$ws.Range("B94").Value = "RDM, petits États insulaires en développement"
$ws.Range("B96").Value = "RDM, pays en développement sans littoral"
$ws.Range("B97").Value = "Afrique, États fragiles"
$ws.Range("B98").Value = "RDM, États fragiles"

# --- 3: swap the two footnote paragraphs -----------------------------------
$ws.Range("A104").Value = "Les exportations / importations des cinq pays membres de l'Union douanière de l'Afrique australe (Afrique du Sud, Botswana, Eswatini, Lesotho et Namibie) sont déclarées en tant qu'exportations / importations pour l'Afrique du Sud dans ces données afin d'améliorer la cohérence des données."
$ws.Range("A105").Value = "Responsabilité : Ce document, ainsi que les données et cartes qu'il peut comprendre, sont sans préjudice du statut de tout territoire, de la souveraineté s'exerçant sur ce dernier, du tracé des frontières et limites internationales, et du nom de tout territoire, ville ou région."

# --- 4: refreshed data for row 92 ("RDM, pays les moins avancés") ---------
$ws.Range("C92").Value = 17394.689651000001
$ws.Range("D92").Value = 13759.125681
$ws.Range("E92").Value = 79731.397599999997
$ws.Range("F92").Value = 110885.21293199999
$ws.Range("G92").Value = 14130.729542999999
$ws.Range("H92").Value = 96990.421031000005
$ws.Range("I92").Value = 41725.939704999997
$ws.Range("J92").Value = 152847.090279

# --- 5: saved window size (best effort; cosmetic view-state) --------------
$win = $wb.Windows.Item(1)
$win.Width = 960
$win.Height = 279.5
